$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark Milestone "I" (and Milestone Complete "X") for the newly achieved
# rubric items: lighting (rows 30-33), a model loaded from file (row 5),
# and a texture (row 18).
$rows = @(5, 18, 30, 31, 32, 33)
foreach ($r in $rows) {
    $ws.Range("E$r").Value = "I"
    $ws.Range("F$r").Value = "X"
}

# Restore the selection to where the author was last working.
$ws.Range("E91").Select()
